# "Class updated to 2024" - Dropout slide wording tweak:
# "Can be applied (or less) independently to each module"
#   -> "Can be applied (or not) independently to each module"
#
# PowerPoint splits the edited fragment into its own run (the untouched
# text on either side keeps its original run/formatting), so we locate
# the sentence, then re-type just the "or less) " -> "or not) " portion
# through a Characters() sub-range rather than rewriting the whole
# paragraph (which would collapse it back into a single run).

$p = $ppt.ActivePresentation

$needle = "Can be applied (or less) independently to each module"
$oldFrag = "or less) "
$newFrag = "or not) "

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if (-not $shape.HasTextFrame) { continue }

        $tr = $shape.TextFrame.TextRange
        $full = $tr.Text
        $pos = $full.IndexOf($needle)
        if ($pos -ge 0) {
            $fragPos = $pos + $needle.IndexOf($oldFrag)
            $frag = $tr.Characters($fragPos + 1, $oldFrag.Length)
            $frag.Text = $newFrag
        }
    }
}
